$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (Tucker Carlson article) had a stale relative-time placeholder
# ("11h ago") in the date column instead of the proper constant date
# value. Replace it with the real scraped date "May 10", matching the
# other rows' date format.
$ws.Range("B6").Value = "May 10"
